# Añadido el reloj: rellenar las 3 celdas vacías de la fila 1 de la tabla
# del reloj (slide 1) con ":" "s" "s" para formar "hh:mm:ss".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

$tbl.Cell(1,6).Shape.TextFrame.TextRange.Text = ":"
$tbl.Cell(1,7).Shape.TextFrame.TextRange.Text = "s"
$tbl.Cell(1,8).Shape.TextFrame.TextRange.Text = "s"
